$wb = $excel.ActiveWorkbook

# ===== Sheet: Overview =====
$ws = $wb.Worksheets.Item("Overview")
$ws.Cells.Item(5,1).Value = "'5aae2a01-bf15-4f6b-bfed-94e1b815c25a.md"
$ws.Cells.Item(5,2).Value = "'e2e\5aae2a01-bf15-4f6b-bfed-94e1b815c25a.md"
$ws.Cells.Item(5,3).Value = "'.md"
$ws.Cells.Item(5,4).Value = "'"
$ws.Cells.Item(5,5).Value = "'Ready for handoff"
$ws.Cells.Item(5,6).Value = "'Ready for handoff"
$ws.Cells.Item(5,7).Value = "'2016-08-21 02:46:08"
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b6a6f1c5c6b4a3e8f0d9c7b1a5e3d2f4c6b8a9d/e2e/5aae2a01-bf15-4f6b-bfed-94e1b815c25a.md", "", "", "e2e\5aae2a01-bf15-4f6b-bfed-94e1b815c25a.md")
$ws.Cells.Item(6,1).Value = "'63f887fa-7ab2-4f2e-b269-e28d0aa51012.md"
$ws.Cells.Item(6,2).Value = "'e2e\63f887fa-7ab2-4f2e-b269-e28d0aa51012.md"
$ws.Cells.Item(6,3).Value = "'.md"
$ws.Cells.Item(6,4).Value = "'"
$ws.Cells.Item(6,5).Value = "'Ready for handoff"
$ws.Cells.Item(6,6).Value = "'Ready for handoff"
$ws.Cells.Item(6,7).Value = "'2016-08-21 02:44:32"
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c834495c83317b2e49b255d9bfe7c2115f74f45/e2e/63f887fa-7ab2-4f2e-b269-e28d0aa51012.md", "", "", "e2e\63f887fa-7ab2-4f2e-b269-e28d0aa51012.md")
$ws.Cells.Item(7,1).Value = "'fc3bd347-9189-4383-a242-abca36ae7e57.md"
$ws.Cells.Item(7,2).Value = "'e2e\fc3bd347-9189-4383-a242-abca36ae7e57.md"
$ws.Cells.Item(7,3).Value = "'.md"
$ws.Cells.Item(7,4).Value = "'"
$ws.Cells.Item(7,5).Value = "'Ready for handoff"
$ws.Cells.Item(7,6).Value = "'Ready for handoff"
$ws.Cells.Item(7,7).Value = "'2016-08-21 02:46:08"
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e4d2c8b9a1f3e5d6c0b4a8f2e1d9c7b5a3f6e8d/e2e/fc3bd347-9189-4383-a242-abca36ae7e57.md", "", "", "e2e\fc3bd347-9189-4383-a242-abca36ae7e57.md")
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G7"))

# ===== Sheet: zh-cn =====
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Cells.Item(5,1).Value = "'5aae2a01-bf15-4f6b-bfed-94e1b815c25a.md"
$ws.Cells.Item(5,2).Value = "'.md"
$ws.Cells.Item(5,3).Value = "'Ready for handoff"
$ws.Cells.Item(5,4).Value = "'e2e"
$ws.Cells.Item(5,5).Value = "'ht"
$ws.Cells.Item(5,6).Value = "'False"
$ws.Cells.Item(5,7).Value = "'5aae2a01-bf15-4f6b-bfed-94e1b815c25a.6abeb009ca820be9e6c616527eef2c6f15c3f337.zh-cn.xlf"
$ws.Cells.Item(5,8).Value = "'2016-08-21 02:45:59"
$ws.Cells.Item(5,9).Value = "'"
$ws.Cells.Item(5,10).Value = "'"
$ws.Cells.Item(5,11).Value = "'0001-01-01 00:00:00"
$ws.Cells.Item(5,12).Value = "'"
$ws.Cells.Item(5,13).Value = "'True"
$ws.Cells.Item(5,14).Value = "'"
$ws.Cells.Item(5,15).Value = "'False"
$ws.Cells.Item(5,16).Value = "'"
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b6a6f1c5c6b4a3e8f0d9c7b1a5e3d2f4c6b8a9d/e2e/5aae2a01-bf15-4f6b-bfed-94e1b815c25a.md", "", "", "5aae2a01-bf15-4f6b-bfed-94e1b815c25a.md")
$ws.Cells.Item(6,1).Value = "'63f887fa-7ab2-4f2e-b269-e28d0aa51012.md"
$ws.Cells.Item(6,2).Value = "'.md"
$ws.Cells.Item(6,3).Value = "'Ready for handoff"
$ws.Cells.Item(6,4).Value = "'e2e"
$ws.Cells.Item(6,5).Value = "'ht"
$ws.Cells.Item(6,6).Value = "'False"
$ws.Cells.Item(6,7).Value = "'63f887fa-7ab2-4f2e-b269-e28d0aa51012.a8095180b0c8425b5d093193c22e36115d5319e1.zh-cn.xlf"
$ws.Cells.Item(6,8).Value = "'2016-08-21 02:44:28"
$ws.Cells.Item(6,9).Value = "'"
$ws.Cells.Item(6,10).Value = "'"
$ws.Cells.Item(6,11).Value = "'0001-01-01 00:00:00"
$ws.Cells.Item(6,12).Value = "'"
$ws.Cells.Item(6,13).Value = "'True"
$ws.Cells.Item(6,14).Value = "'"
$ws.Cells.Item(6,15).Value = "'False"
$ws.Cells.Item(6,16).Value = "'"
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c834495c83317b2e49b255d9bfe7c2115f74f45/e2e/63f887fa-7ab2-4f2e-b269-e28d0aa51012.md", "", "", "63f887fa-7ab2-4f2e-b269-e28d0aa51012.md")
$ws.Cells.Item(7,1).Value = "'fc3bd347-9189-4383-a242-abca36ae7e57.md"
$ws.Cells.Item(7,2).Value = "'.md"
$ws.Cells.Item(7,3).Value = "'Ready for handoff"
$ws.Cells.Item(7,4).Value = "'e2e"
$ws.Cells.Item(7,5).Value = "'ht"
$ws.Cells.Item(7,6).Value = "'False"
$ws.Cells.Item(7,7).Value = "'fc3bd347-9189-4383-a242-abca36ae7e57.5ca95beaaf8809d0a37c82df808f175f56d729cb.zh-cn.xlf"
$ws.Cells.Item(7,8).Value = "'2016-08-21 02:45:59"
$ws.Cells.Item(7,9).Value = "'"
$ws.Cells.Item(7,10).Value = "'"
$ws.Cells.Item(7,11).Value = "'0001-01-01 00:00:00"
$ws.Cells.Item(7,12).Value = "'"
$ws.Cells.Item(7,13).Value = "'True"
$ws.Cells.Item(7,14).Value = "'"
$ws.Cells.Item(7,15).Value = "'False"
$ws.Cells.Item(7,16).Value = "'"
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e4d2c8b9a1f3e5d6c0b4a8f2e1d9c7b5a3f6e8d/e2e/fc3bd347-9189-4383-a242-abca36ae7e57.md", "", "", "fc3bd347-9189-4383-a242-abca36ae7e57.md")
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P7"))

# ===== Sheet: de-de =====
$ws = $wb.Worksheets.Item("de-de")
$ws.Cells.Item(5,1).Value = "'5aae2a01-bf15-4f6b-bfed-94e1b815c25a.md"
$ws.Cells.Item(5,2).Value = "'.md"
$ws.Cells.Item(5,3).Value = "'Ready for handoff"
$ws.Cells.Item(5,4).Value = "'e2e"
$ws.Cells.Item(5,5).Value = "'ht"
$ws.Cells.Item(5,6).Value = "'False"
$ws.Cells.Item(5,7).Value = "'5aae2a01-bf15-4f6b-bfed-94e1b815c25a.6abeb009ca820be9e6c616527eef2c6f15c3f337.de-de.xlf"
$ws.Cells.Item(5,8).Value = "'2016-08-21 02:46:08"
$ws.Cells.Item(5,9).Value = "'"
$ws.Cells.Item(5,10).Value = "'"
$ws.Cells.Item(5,11).Value = "'0001-01-01 00:00:00"
$ws.Cells.Item(5,12).Value = "'"
$ws.Cells.Item(5,13).Value = "'True"
$ws.Cells.Item(5,14).Value = "'"
$ws.Cells.Item(5,15).Value = "'False"
$ws.Cells.Item(5,16).Value = "'"
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b6a6f1c5c6b4a3e8f0d9c7b1a5e3d2f4c6b8a9d/e2e/5aae2a01-bf15-4f6b-bfed-94e1b815c25a.md", "", "", "5aae2a01-bf15-4f6b-bfed-94e1b815c25a.md")
$ws.Cells.Item(6,1).Value = "'63f887fa-7ab2-4f2e-b269-e28d0aa51012.md"
$ws.Cells.Item(6,2).Value = "'.md"
$ws.Cells.Item(6,3).Value = "'Ready for handoff"
$ws.Cells.Item(6,4).Value = "'e2e"
$ws.Cells.Item(6,5).Value = "'ht"
$ws.Cells.Item(6,6).Value = "'False"
$ws.Cells.Item(6,7).Value = "'63f887fa-7ab2-4f2e-b269-e28d0aa51012.a8095180b0c8425b5d093193c22e36115d5319e1.de-de.xlf"
$ws.Cells.Item(6,8).Value = "'2016-08-21 02:44:32"
$ws.Cells.Item(6,9).Value = "'"
$ws.Cells.Item(6,10).Value = "'"
$ws.Cells.Item(6,11).Value = "'0001-01-01 00:00:00"
$ws.Cells.Item(6,12).Value = "'"
$ws.Cells.Item(6,13).Value = "'True"
$ws.Cells.Item(6,14).Value = "'"
$ws.Cells.Item(6,15).Value = "'False"
$ws.Cells.Item(6,16).Value = "'"
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c834495c83317b2e49b255d9bfe7c2115f74f45/e2e/63f887fa-7ab2-4f2e-b269-e28d0aa51012.md", "", "", "63f887fa-7ab2-4f2e-b269-e28d0aa51012.md")
$ws.Cells.Item(7,1).Value = "'fc3bd347-9189-4383-a242-abca36ae7e57.md"
$ws.Cells.Item(7,2).Value = "'.md"
$ws.Cells.Item(7,3).Value = "'Ready for handoff"
$ws.Cells.Item(7,4).Value = "'e2e"
$ws.Cells.Item(7,5).Value = "'ht"
$ws.Cells.Item(7,6).Value = "'False"
$ws.Cells.Item(7,7).Value = "'fc3bd347-9189-4383-a242-abca36ae7e57.5ca95beaaf8809d0a37c82df808f175f56d729cb.de-de.xlf"
$ws.Cells.Item(7,8).Value = "'2016-08-21 02:46:08"
$ws.Cells.Item(7,9).Value = "'"
$ws.Cells.Item(7,10).Value = "'"
$ws.Cells.Item(7,11).Value = "'0001-01-01 00:00:00"
$ws.Cells.Item(7,12).Value = "'"
$ws.Cells.Item(7,13).Value = "'True"
$ws.Cells.Item(7,14).Value = "'"
$ws.Cells.Item(7,15).Value = "'False"
$ws.Cells.Item(7,16).Value = "'"
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e4d2c8b9a1f3e5d6c0b4a8f2e1d9c7b5a3f6e8d/e2e/fc3bd347-9189-4383-a242-abca36ae7e57.md", "", "", "fc3bd347-9189-4383-a242-abca36ae7e57.md")
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P7"))
